# Update "想去人数" (want-to-go count) figures in column F across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1161
$ws1.Range("F8").Value = 1072
$ws1.Range("F9").Value = 1654
$ws1.Range("F12").Value = 1784
$ws1.Range("F13").Value = 462
$ws1.Range("F16").Value = 6390
$ws1.Range("F17").Value = 6390
$ws1.Range("F28").Value = 738
$ws1.Range("F29").Value = 289
$ws1.Range("F34").Value = 3877

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 400
$ws2.Range("F12").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1161
$ws4.Range("F13").Value = 1654
$ws4.Range("F16").Value = 1784
$ws4.Range("F18").Value = 462
$ws4.Range("F21").Value = 6
$ws4.Range("F22").Value = 6390
$ws4.Range("F23").Value = 6390
$ws4.Range("F34").Value = 738
$ws4.Range("F36").Value = 289
$ws4.Range("F45").Value = 3877
